$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save row 15 values for the columns that change
$a15 = $ws.Range("A15").Value2
$b15 = $ws.Range("B15").Value2
$d15 = $ws.Range("D15").Value2
$e15 = $ws.Range("E15").Value2
$f15 = $ws.Range("F15").Value2
$g15 = $ws.Range("G15").Value2
$h15 = $ws.Range("H15").Value2
$q15 = $ws.Range("Q15").Value2
$r15 = $ws.Range("R15").Value2

# Save row 16 values
$a16 = $ws.Range("A16").Value2
$b16 = $ws.Range("B16").Value2
$d16 = $ws.Range("D16").Value2
$e16 = $ws.Range("E16").Value2
$f16 = $ws.Range("F16").Value2
$g16 = $ws.Range("G16").Value2
$h16 = $ws.Range("H16").Value2
$q16 = $ws.Range("Q16").Value2
$r16 = $ws.Range("R16").Value2

# Swap: row15 gets row16's original data, row16 gets row15's original data
$ws.Range("A15").Value2 = $a16
$ws.Range("B15").Value2 = $b16
$ws.Range("D15").Value2 = $d16
$ws.Range("E15").Value2 = $e16
$ws.Range("F15").Value2 = $f16
$ws.Range("G15").Value2 = $g16
$ws.Range("H15").Value2 = $h16
$ws.Range("Q15").Value2 = $q16
$ws.Range("R15").Value2 = $r16

$ws.Range("A16").Value2 = $a15
$ws.Range("B16").Value2 = $b15
$ws.Range("D16").Value2 = $d15
$ws.Range("E16").Value2 = $e15
$ws.Range("F16").Value2 = $f15
$ws.Range("G16").Value2 = $g15
$ws.Range("H16").Value2 = $h15
$ws.Range("Q16").Value2 = $q15
$ws.Range("R16").Value2 = $r15
